$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings like "1.00" or "0.0000240" are not coerced into numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('B2').Formula = 'Bitcoin'
$ws.Range('C2').Formula = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Formula = '67.165.81'
$ws.Range('E2').Formula = '  -3.42%  '

$ws.Range('B3').Formula = 'Ethereum'
$ws.Range('C3').Formula = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Formula = '3.681.53'
$ws.Range('E3').Formula = '  -2.56%  '

$ws.Range('B4').Formula = 'TetherUSD'
$ws.Range('C4').Formula = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Formula = '1.00'
$ws.Range('E4').Formula = '  +0.18%  '

$ws.Range('B5').Formula = 'BNB'
$ws.Range('C5').Formula = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Formula = '590.90'
$ws.Range('E5').Formula = '  -4.00%  '

$ws.Range('B6').Formula = 'Solana'
$ws.Range('C6').Formula = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Formula = '166.43'
$ws.Range('E6').Formula = '  -6.21%  '

$ws.Range('B7').Formula = 'LidoStakedEther'
$ws.Range('C7').Formula = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Formula = '3.680.95'
$ws.Range('E7').Formula = '  -2.52%  '

$ws.Range('B8').Formula = 'USDC'
$ws.Range('C8').Formula = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Formula = '1.00'
$ws.Range('E8').Formula = '  +0.01%  '

$ws.Range('B9').Formula = 'XRP'
$ws.Range('C9').Formula = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').Formula = '0.521'
$ws.Range('E9').Formula = '  -0.76%  '

$ws.Range('B10').Formula = 'Dogecoin'
$ws.Range('C10').Formula = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Formula = '0.161'
$ws.Range('E10').Formula = '  -2.74%  '

$ws.Range('B11').Formula = 'Toncoin'
$ws.Range('C11').Formula = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Formula = '6.12'
$ws.Range('E11').Formula = '  -5.24%  '

$ws.Range('B12').Formula = 'Cardano'
$ws.Range('C12').Formula = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D12').Formula = '0.458'
$ws.Range('E12').Formula = '  -5.31%  '

$ws.Range('B13').Formula = 'Avalanche'
$ws.Range('C13').Formula = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').Formula = '37.68'
$ws.Range('E13').Formula = '  -5.27%  '

$ws.Range('B14').Formula = 'ShibaInu'
$ws.Range('C14').Formula = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Formula = '0.0000240'
$ws.Range('E14').Formula = '  -5.71%  '

$ws.Range('B15').Formula = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Formula = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Formula = '4.312.65'
$ws.Range('E15').Formula = '  -1.93%  '

$ws.Range('B16').Formula = 'WrappedEther'
$ws.Range('C16').Formula = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Formula = '3.697.76'
$ws.Range('E16').Formula = '  -1.82%  '

$ws.Range('B17').Formula = 'WrappedBTC'
$ws.Range('C17').Formula = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Formula = '67.245.16'
$ws.Range('E17').Formula = '  -3.34%  '

$ws.Range('B18').Formula = 'TRON'
$ws.Range('C18').Formula = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Formula = '0.114'
$ws.Range('E18').Formula = '  -4.04%  '

$ws.Range('B19').Formula = 'Polkadot'
$ws.Range('C19').Formula = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D19').Formula = '7.08'
$ws.Range('E19').Formula = '  -6.41%  '

$ws.Range('B20').Formula = 'Chainlink'
$ws.Range('C20').Formula = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Formula = '16.99'
$ws.Range('E20').Formula = '  +2.30%  '

$ws.Range('B21').Formula = 'BitcoinCash'
$ws.Range('C21').Formula = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Formula = '486.60'
$ws.Range('E21').Formula = '  -4.55%  '

$ws.Range('B22').Formula = 'Uniswap'
$ws.Range('C22').Formula = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Formula = '9.07'
$ws.Range('E22').Formula = '  -4.68%  '

$ws.Range('B23').Formula = 'Polygon'
$ws.Range('C23').Formula = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Formula = '0.719'
$ws.Range('E23').Formula = '  -2.11%  '

$ws.Range('B24').Formula = 'Litecoin'
$ws.Range('C24').Formula = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').Formula = '84.83'
$ws.Range('E24').Formula = '  -1.69%  '

$ws.Range('B25').Formula = 'Fetch.AI'
$ws.Range('C25').Formula = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').Formula = '2.30'
$ws.Range('E25').Formula = '  -7.23%  '

$ws.Range('B26').Formula = 'PEPE'
$ws.Range('C26').Formula = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Formula = '0.0000140'
$ws.Range('E26').Formula = '  -0.45%  '

$ws.Range('B27').Formula = 'InternetComputer(DFINITY)'
$ws.Range('C27').Formula = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Formula = '12.09'
$ws.Range('E27').Formula = '  -6.26%  '

$ws.Range('B28').Formula = 'Dai'
$ws.Range('C28').Formula = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').Formula = '0.996'
$ws.Range('E28').Formula = '  -0.41%  '

$ws.Range('B29').Formula = 'RenderToken'
$ws.Range('C29').Formula = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').Formula = '9.91'
$ws.Range('E29').Formula = '  -6.19%  '

$ws.Range('B30').Formula = 'PancakeSwap'
$ws.Range('C30').Formula = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Formula = '2.91'
$ws.Range('E30').Formula = '  -3.34%  '

$ws.Range('B31').Formula = 'ImmutableX'
$ws.Range('C31').Formula = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Formula = '2.35'
$ws.Range('E31').Formula = '  -6.67%  '

$ws.Range('B32').Formula = 'NEARProtocol'
$ws.Range('C32').Formula = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Formula = '7.67'
$ws.Range('E32').Formula = '  -5.80%  '

$ws.Range('B33').Formula = 'EthereumClassic'
$ws.Range('C33').Formula = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').Formula = '31.55'
$ws.Range('E33').Formula = '  +1.28%  '

$ws.Range('B34').Formula = 'WrappedeETH'
$ws.Range('C34').Formula = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D34').Formula = '3.833.94'
$ws.Range('E34').Formula = '  -1.96%  '

$ws.Range('B35').Formula = 'RenzoRestakedETH'
$ws.Range('C35').Formula = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Formula = '3.628.63'
$ws.Range('E35').Formula = '  -2.16%  '

$ws.Range('B36').Formula = 'Hedera'
$ws.Range('C36').Formula = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Formula = '0.107'
$ws.Range('E36').Formula = '  -7.23%  '

$ws.Range('B37').Formula = 'FirstDigitalUSD'
$ws.Range('C37').Formula = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').Formula = '1.00'
$ws.Range('E37').Formula = '  +0.16%  '

$ws.Range('B38').Formula = 'Mantle'
$ws.Range('C38').Formula = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').Formula = '0.995'
$ws.Range('E38').Formula = '  -5.18%  '

$ws.Range('B39').Formula = 'Filecoin'
$ws.Range('C39').Formula = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').Formula = '5.72'
$ws.Range('E39').Formula = '  -6.85%  '

$ws.Range('B40').Formula = 'Kaspa'
$ws.Range('C40').Formula = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Formula = '0.131'
$ws.Range('E40').Formula = '  -7.99%  '

$ws.Range('B41').Formula = 'TheGraph'
$ws.Range('C41').Formula = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Formula = '0.320'
$ws.Range('E41').Formula = '  -5.86%  '

$ws.Range('B42').Formula = 'Bittensor'
$ws.Range('C42').Formula = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').Formula = '444.26'
$ws.Range('E42').Formula = '  -6.77%  '

$ws.Range('B43').Formula = 'OKB'
$ws.Range('C43').Formula = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Formula = '48.79'
$ws.Range('E43').Formula = '  -1.95%  '

$ws.Range('B44').Formula = 'Stacks'
$ws.Range('C44').Formula = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D44').Formula = '1.95'
$ws.Range('E44').Formula = '  -5.87%  '

$ws.Range('B45').Formula = 'dogwifhat'
$ws.Range('C45').Formula = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Formula = '2.77'
$ws.Range('E45').Formula = '  -7.51%  '

$ws.Range('B46').Formula = 'Cosmos'
$ws.Range('C46').Formula = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').Formula = '8.25'
$ws.Range('E46').Formula = '  -3.92%  '

$ws.Range('B47').Formula = 'USDe'
$ws.Range('C47').Formula = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Formula = '1.00'
$ws.Range('E47').Formula = '  +0.05%  '

$ws.Range('B48').Formula = 'Arweave'
$ws.Range('C48').Formula = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D48').Formula = '39.78'
$ws.Range('E48').Formula = '  -10.08%  '

$ws.Range('B49').Formula = 'Monero'
$ws.Range('C49').Formula = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D49').Formula = '140.58'
$ws.Range('E49').Formula = '  +0.71%  '

$ws.Range('B50').Formula = 'Maker'
$ws.Range('C50').Formula = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Formula = '2.770.93'
$ws.Range('E50').Formula = '  -5.96%  '

$ws.Range('B51').Formula = 'VeChain'
$ws.Range('C51').Formula = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Formula = '0.0344'
$ws.Range('E51').Formula = '  -4.87%  '
